# Apply the Flashscore weekly-odds refresh for 2024-11-04.
# A new fixture (Corinthians vs Palmeiras) is inserted logically at row 6;
# the fixtures that used to occupy rows 6-8 move down to rows 7-9 (their odds
# are refreshed too), the old row-9 fixture (Hapoel Haifa vs Maccabi Haifa) is
# dropped, and several other rows get updated odds values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet


# Row 3 (rXmhxCKr - Lazio vs Cagliari): odds refresh
$ws.Cells.Item(3, 7).Value = 1.55
$ws.Cells.Item(3, 8).Value = 3.9
$ws.Cells.Item(3, 9).Value = 6.25
$ws.Cells.Item(3, 10).Value = 2.1
$ws.Cells.Item(3, 12).Value = 5.5
$ws.Cells.Item(3, 13).Value = 1.05
$ws.Cells.Item(3, 14).Value = 11
$ws.Cells.Item(3, 15).Value = 1.25
$ws.Cells.Item(3, 16).Value = 4
$ws.Cells.Item(3, 17).Value = 1.8
$ws.Cells.Item(3, 18).Value = 2
$ws.Cells.Item(3, 21).Value = 1.8
$ws.Cells.Item(3, 22).Value = 1.95
$ws.Cells.Item(3, 34).Value = 17
$ws.Cells.Item(3, 36).Value = 19
$ws.Cells.Item(3, 38).Value = 41
$ws.Cells.Item(3, 39).Value = 41
$ws.Cells.Item(3, 43).Value = 23
$ws.Cells.Item(3, 49).Value = 126
$ws.Cells.Item(3, 51).Value = 29
$ws.Cells.Item(3, 53).Value = 101

# Row 4 (zeCiJHet - Celta Vigo vs Getafe): odds refresh
$ws.Cells.Item(4, 7).Value = 2.15
$ws.Cells.Item(4, 8).Value = 3
$ws.Cells.Item(4, 9).Value = 3.8
$ws.Cells.Item(4, 10).Value = 3
$ws.Cells.Item(4, 12).Value = 4.5
$ws.Cells.Item(4, 15).Value = 1.53
$ws.Cells.Item(4, 16).Value = 2.5
$ws.Cells.Item(4, 25).Value = 10
$ws.Cells.Item(4, 34).Value = 8
$ws.Cells.Item(4, 35).Value = 17
$ws.Cells.Item(4, 40).Value = 4
$ws.Cells.Item(4, 41).Value = 13
$ws.Cells.Item(4, 55).Value = 351

# Row 6 (new fixture fT8rSK5A - Corinthians vs Palmeiras)
$ws.Cells.Item(6, 1).Value = 'fT8rSK5A'
$ws.Cells.Item(6, 3).Value = '20:00'
$ws.Cells.Item(6, 4).Value = 'BRAZIL - SERIE A BETANO'
$ws.Cells.Item(6, 5).Value = 'Corinthians'
$ws.Cells.Item(6, 6).Value = 'Palmeiras'
$ws.Cells.Item(6, 7).Value = 3.1
$ws.Cells.Item(6, 8).Value = 3.25
$ws.Cells.Item(6, 9).Value = 2.38
$ws.Cells.Item(6, 10).Value = 3.75
$ws.Cells.Item(6, 11).Value = 2.05
$ws.Cells.Item(6, 13).Value = 1.08
$ws.Cells.Item(6, 14).Value = 8
$ws.Cells.Item(6, 15).Value = 1.4
$ws.Cells.Item(6, 16).Value = 3
$ws.Cells.Item(6, 17).Value = 2.15
$ws.Cells.Item(6, 18).Value = 1.67
$ws.Cells.Item(6, 19).Value = 1.44
$ws.Cells.Item(6, 20).Value = 2.63
$ws.Cells.Item(6, 21).Value = 1.91
$ws.Cells.Item(6, 22).Value = 1.91
$ws.Cells.Item(6, 23).Value = 8.5
$ws.Cells.Item(6, 25).Value = 11
$ws.Cells.Item(6, 26).Value = 34
$ws.Cells.Item(6, 27).Value = 26
$ws.Cells.Item(6, 28).Value = 34
$ws.Cells.Item(6, 29).Value = 8.5
$ws.Cells.Item(6, 31).Value = 15
$ws.Cells.Item(6, 32).Value = 51
$ws.Cells.Item(6, 33).Value = 301
$ws.Cells.Item(6, 34).Value = 7.5
$ws.Cells.Item(6, 35).Value = 11
$ws.Cells.Item(6, 36).Value = 9.5
$ws.Cells.Item(6, 37).Value = 23
$ws.Cells.Item(6, 38).Value = 21
$ws.Cells.Item(6, 39).Value = 34
$ws.Cells.Item(6, 41).Value = 17
$ws.Cells.Item(6, 42).Value = 29
$ws.Cells.Item(6, 43).Value = 51
$ws.Cells.Item(6, 44).Value = 81
$ws.Cells.Item(6, 45).Value = 201
$ws.Cells.Item(6, 46).Value = 2.63
$ws.Cells.Item(6, 47).Value = 8
$ws.Cells.Item(6, 48).Value = 51
$ws.Cells.Item(6, 49).Value = 126
$ws.Cells.Item(6, 50).Value = 4.33
$ws.Cells.Item(6, 52).Value = 23
$ws.Cells.Item(6, 53).Value = 41
$ws.Cells.Item(6, 54).Value = 67
$ws.Cells.Item(6, 55).Value = 201
$ws.Cells.Item(6, 56).Value = 126

# Row 7 (tGmmLIvL - Operario vs Sport Recife, formerly row 6)
$ws.Cells.Item(7, 1).Value = 'tGmmLIvL'
$ws.Cells.Item(7, 3).Value = '19:00'
$ws.Cells.Item(7, 4).Value = 'BRAZIL - SERIE B'
$ws.Cells.Item(7, 5).Value = 'Operario'
$ws.Cells.Item(7, 6).Value = 'Sport Recife'
$ws.Cells.Item(7, 7).Value = 3.4
$ws.Cells.Item(7, 8).Value = 3
$ws.Cells.Item(7, 9).Value = 2.25
$ws.Cells.Item(7, 10).Value = 4.33
$ws.Cells.Item(7, 11).Value = 1.91
$ws.Cells.Item(7, 12).Value = 3.1
$ws.Cells.Item(7, 13).Value = 1.1
$ws.Cells.Item(7, 14).Value = 7
$ws.Cells.Item(7, 15).Value = 1.5
$ws.Cells.Item(7, 16).Value = 2.5
$ws.Cells.Item(7, 17).Value = 2.6
$ws.Cells.Item(7, 18).Value = 1.48
$ws.Cells.Item(7, 19).Value = 1.57
$ws.Cells.Item(7, 20).Value = 2.25
$ws.Cells.Item(7, 21).Value = 2.2
$ws.Cells.Item(7, 22).Value = 1.62
$ws.Cells.Item(7, 23).Value = 7.5
$ws.Cells.Item(7, 25).Value = 13
$ws.Cells.Item(7, 26).Value = 41
$ws.Cells.Item(7, 27).Value = 34
$ws.Cells.Item(7, 29).Value = 6.5
$ws.Cells.Item(7, 30).Value = 6
$ws.Cells.Item(7, 31).Value = 19
$ws.Cells.Item(7, 32).Value = 67
$ws.Cells.Item(7, 34).Value = 6
$ws.Cells.Item(7, 35).Value = 9.5
$ws.Cells.Item(7, 36).Value = 10
$ws.Cells.Item(7, 38).Value = 23
$ws.Cells.Item(7, 39).Value = 41
$ws.Cells.Item(7, 41).Value = 21
$ws.Cells.Item(7, 42).Value = 34
$ws.Cells.Item(7, 44).Value = 126
$ws.Cells.Item(7, 45).Value = 351
$ws.Cells.Item(7, 46).Value = 2.25
$ws.Cells.Item(7, 47).Value = 9
$ws.Cells.Item(7, 48).Value = 67
$ws.Cells.Item(7, 49).Value = 81
$ws.Cells.Item(7, 50).Value = 4
$ws.Cells.Item(7, 52).Value = 29
$ws.Cells.Item(7, 53).Value = 51
$ws.Cells.Item(7, 54).Value = 81
$ws.Cells.Item(7, 55).Value = 251
$ws.Cells.Item(7, 56).Value = 81

# Row 8 (zuChyeRl - Tecnico U. vs U. Catolica, formerly row 7)
$ws.Cells.Item(8, 1).Value = 'zuChyeRl'
$ws.Cells.Item(8, 3).Value = '17:00'
$ws.Cells.Item(8, 4).Value = 'ECUADOR - LIGA PRO'
$ws.Cells.Item(8, 5).Value = 'Tecnico U.'
$ws.Cells.Item(8, 6).Value = 'U. Catolica'
$ws.Cells.Item(8, 7).Value = 3.2
$ws.Cells.Item(8, 8).Value = 3.3
$ws.Cells.Item(8, 9).Value = 2.2
$ws.Cells.Item(8, 10).Value = 4
$ws.Cells.Item(8, 11).Value = 2.05
$ws.Cells.Item(8, 12).Value = 3
$ws.Cells.Item(8, 13).Value = 1.07
$ws.Cells.Item(8, 14).Value = 9
$ws.Cells.Item(8, 15).Value = 1.36
$ws.Cells.Item(8, 16).Value = 3
$ws.Cells.Item(8, 17).Value = 2.15
$ws.Cells.Item(8, 18).Value = 1.67
$ws.Cells.Item(8, 21).Value = 1.91
$ws.Cells.Item(8, 22).Value = 1.8
$ws.Cells.Item(8, 23).Value = 8.5
$ws.Cells.Item(8, 24).Value = 15
$ws.Cells.Item(8, 25).Value = 12
$ws.Cells.Item(8, 26).Value = 34
$ws.Cells.Item(8, 27).Value = 29
$ws.Cells.Item(8, 28).Value = 41
$ws.Cells.Item(8, 29).Value = 8.5
$ws.Cells.Item(8, 31).Value = 17
$ws.Cells.Item(8, 32).Value = 51
$ws.Cells.Item(8, 33).Value = 351
$ws.Cells.Item(8, 34).Value = 7
$ws.Cells.Item(8, 35).Value = 10
$ws.Cells.Item(8, 36).Value = 9.5
$ws.Cells.Item(8, 37).Value = 21
$ws.Cells.Item(8, 40).Value = 5
$ws.Cells.Item(8, 41).Value = 19
$ws.Cells.Item(8, 42).Value = 29
$ws.Cells.Item(8, 43).Value = 67
$ws.Cells.Item(8, 44).Value = 101
$ws.Cells.Item(8, 45).Value = 251
$ws.Cells.Item(8, 47).Value = 8.5
$ws.Cells.Item(8, 48).Value = 51
$ws.Cells.Item(8, 49).Value = 51
$ws.Cells.Item(8, 50).Value = 4.33
$ws.Cells.Item(8, 51).Value = 13
$ws.Cells.Item(8, 52).Value = 23
$ws.Cells.Item(8, 56).Value = 51

# Row 9 (4WNXteh2 - AC Ajaccio vs Metz, formerly row 8)
$ws.Cells.Item(9, 1).Value = '4WNXteh2'
$ws.Cells.Item(9, 3).Value = '16:45'
$ws.Cells.Item(9, 4).Value = 'FRANCE - LIGUE 2'
$ws.Cells.Item(9, 5).Value = 'AC Ajaccio'
$ws.Cells.Item(9, 6).Value = 'Metz'
$ws.Cells.Item(9, 7).Value = 4.33
$ws.Cells.Item(9, 8).Value = 3.25
$ws.Cells.Item(9, 9).Value = 1.9
$ws.Cells.Item(9, 10).Value = 5
$ws.Cells.Item(9, 11).Value = 2
$ws.Cells.Item(9, 12).Value = 2.63
$ws.Cells.Item(9, 13).Value = 1.08
$ws.Cells.Item(9, 14).Value = 8
$ws.Cells.Item(9, 15).Value = 1.44
$ws.Cells.Item(9, 16).Value = 2.63
$ws.Cells.Item(9, 17).Value = 2.35
$ws.Cells.Item(9, 18).Value = 1.57
$ws.Cells.Item(9, 19).Value = 1.5
$ws.Cells.Item(9, 20).Value = 2.5
$ws.Cells.Item(9, 21).Value = 2.1
$ws.Cells.Item(9, 22).Value = 1.67
$ws.Cells.Item(9, 23).Value = 9.5
$ws.Cells.Item(9, 24).Value = 21
$ws.Cells.Item(9, 25).Value = 15
$ws.Cells.Item(9, 26).Value = 41
$ws.Cells.Item(9, 28).Value = 51
$ws.Cells.Item(9, 29).Value = 7.5
$ws.Cells.Item(9, 30).Value = 6.5
$ws.Cells.Item(9, 31).Value = 19
$ws.Cells.Item(9, 32).Value = 67
$ws.Cells.Item(9, 34).Value = 6
$ws.Cells.Item(9, 35).Value = 8
$ws.Cells.Item(9, 36).Value = 9
$ws.Cells.Item(9, 37).Value = 15
$ws.Cells.Item(9, 38).Value = 19
$ws.Cells.Item(9, 39).Value = 34
$ws.Cells.Item(9, 40).Value = 6
$ws.Cells.Item(9, 41).Value = 26
$ws.Cells.Item(9, 43).Value = 81
$ws.Cells.Item(9, 45).Value = 351
$ws.Cells.Item(9, 46).Value = 2.5
$ws.Cells.Item(9, 47).Value = 9
$ws.Cells.Item(9, 48).Value = 67
$ws.Cells.Item(9, 49).Value = 81
$ws.Cells.Item(9, 50).Value = 3.75
$ws.Cells.Item(9, 51).Value = 11
$ws.Cells.Item(9, 52).Value = 26
$ws.Cells.Item(9, 53).Value = 41
$ws.Cells.Item(9, 54).Value = 67
$ws.Cells.Item(9, 55).Value = 201
$ws.Cells.Item(9, 56).Value = 81

# Row 10 (WrBOmjM7 - Jong PSV vs FC Emmen): odds refresh
$ws.Cells.Item(10, 7).Value = 3.9
$ws.Cells.Item(10, 8).Value = 4
$ws.Cells.Item(10, 14).Value = 19
$ws.Cells.Item(10, 15).Value = 1.14
$ws.Cells.Item(10, 16).Value = 5.5
$ws.Cells.Item(10, 17).Value = 1.5
$ws.Cells.Item(10, 18).Value = 2.5
$ws.Cells.Item(10, 21).Value = 1.5
$ws.Cells.Item(10, 22).Value = 2.5
$ws.Cells.Item(10, 23).Value = 17
$ws.Cells.Item(10, 28).Value = 29
$ws.Cells.Item(10, 29).Value = 19
$ws.Cells.Item(10, 30).Value = 8
$ws.Cells.Item(10, 34).Value = 11
$ws.Cells.Item(10, 35).Value = 11
$ws.Cells.Item(10, 36).Value = 8.5
$ws.Cells.Item(10, 40).Value = 6
$ws.Cells.Item(10, 44).Value = 67
$ws.Cells.Item(10, 49).Value = 301

# Row 11 (zHHFkCie - Jong Utrecht vs Maastricht): odds refresh
$ws.Cells.Item(11, 7).Value = 2.45
$ws.Cells.Item(11, 9).Value = 2.75
$ws.Cells.Item(11, 12).Value = 3.25
$ws.Cells.Item(11, 24).Value = 13
$ws.Cells.Item(11, 26).Value = 23
$ws.Cells.Item(11, 36).Value = 11
$ws.Cells.Item(11, 37).Value = 29
$ws.Cells.Item(11, 38).Value = 21
$ws.Cells.Item(11, 39).Value = 26
$ws.Cells.Item(11, 47).Value = 7.5

# Row 13 (GCW0bGhf - Zaglebie vs Slask Wroclaw): odds refresh
$ws.Cells.Item(13, 17).Value = 2.1
$ws.Cells.Item(13, 18).Value = 1.7

# Row 15 (dQ60rhHO - FC Rapid Bucuresti vs FC Hermannstadt): odds refresh
$ws.Cells.Item(15, 7).Value = 1.53
$ws.Cells.Item(15, 8).Value = 3.7
$ws.Cells.Item(15, 9).Value = 6.25
$ws.Cells.Item(15, 12).Value = 6
$ws.Cells.Item(15, 24).Value = 7
$ws.Cells.Item(15, 26).Value = 11
$ws.Cells.Item(15, 29).Value = 9.5
$ws.Cells.Item(15, 32).Value = 51
$ws.Cells.Item(15, 41).Value = 8
